$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 3.2
$ws.Range("I2").Value = 2.45
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 1.95
$ws.Range("L2").Value = 3.25
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5
$ws.Range("Q2").Value = 1.88
$ws.Range("R2").Value = 1.98
$ws.Range("S2").Value = 2.5
$ws.Range("T2").Value = 1.5
$ws.Range("W2").Value = 5
$ws.Range("X2").Value = 1.17
$ws.Range("AC2").Value = 7.5
$ws.Range("AD2").Value = 15
$ws.Range("AE2").Value = 12
$ws.Range("AF2").Value = 34
$ws.Range("AN2").Value = 6.5
$ws.Range("AO2").Value = 10
$ws.Range("AP2").Value = 10
$ws.Range("AQ2").Value = 23
$ws.Range("AR2").Value = 23

# Row 3 updates
$ws.Range("G3").Value = 2.62
$ws.Range("H3").Value = 3.75
$ws.Range("I3").Value = 2.27
$ws.Range("J3").Value = 3.1
$ws.Range("L3").Value = 2.75
$ws.Range("S3").Value = 1.55
$ws.Range("T3").Value = 2.15
$ws.Range("W3").Value = 2.32
$ws.Range("X3").Value = 1.47
$ws.Range("AA3").Value = 1.5
$ws.Range("AB3").Value = 2.25
$ws.Range("AC3").Value = 12
$ws.Range("AD3").Value = 15.5
$ws.Range("AE3").Value = 10
$ws.Range("AF3").Value = 30
$ws.Range("AG3").Value = 19.5
$ws.Range("AH3").Value = 24
$ws.Range("AI3").Value = 15
$ws.Range("AJ3").Value = 7.6
$ws.Range("AK3").Value = 12
$ws.Range("AL3").Value = 40
$ws.Range("AN3").Value = 11
$ws.Range("AO3").Value = 13
$ws.Range("AP3").Value = 9.25
$ws.Range("AQ3").Value = 23
$ws.Range("AR3").Value = 16.5
